# AB#14631 - Remove the "Status" column from the Immunization Recommendations
# (forecast) table on the Immunizations worksheet.
#
# The forecast table header is in row 8 (Immunization | Due Date | Status)
# with two data rows (9 and 10). The "Status" column (column C) is removed:
#   - C8 ("Status") header text is cleared, but the cell/style stays (it is
#     still part of the styled header row, like the blank G2:I2 cells above).
#   - C9 / C10 (the recommendations[i].status / recommendations[i+1].status
#     template placeholders) are cleared completely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Immunizations")

# Drop the "Status" header label but keep the formatted (empty) cell.
$ws.Range("C8").ClearContents()

# Drop the per-row status placeholders entirely.
$ws.Range("C9").ClearContents()
$ws.Range("C10").ClearContents()

# Move the active selection down to A11, below the forecast table, matching
# the refreshed sheet view saved with the template.
$ws.Range("A11").Select()
